$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Julio de 2020 a las 06:27"

# --- Pakistan (row 15): refreshed case numbers ---
$ws.Range("B15").Value = 255769
$ws.Range("C15").Value = 2165
$ws.Range("D15").Value = 172810
$ws.Range("E15").Value = 77573
$ws.Range("G15").Value = 66
$ws.Range("H15").Value = 5386

# --- Belgica / Kazajistan swap positions (rows 33-34) with refreshed data ---
$ws.Range("A33").Value = "Kazajistan"
$ws.Range("B33").Value = 63514
$ws.Range("C33").Value = 1759
$ws.Range("D33").Value = 38008
$ws.Range("E33").Value = 25131
$ws.Range("H33").Value = 375

$ws.Range("A34").Value = "Belgica"
$ws.Range("B34").Value = 62781
$ws.Range("D34").Value = 17223
$ws.Range("E34").Value = 35771
$ws.Range("H34").Value = 9787

# --- Lesoto / Martinica / Mongolia re-ordered (rows 169-171) with refreshed data ---
$ws.Range("A169").Value = "Mongolia"
$ws.Range("B169").Value = 261
$ws.Range("C169").Value = 18
$ws.Range("D169").Value = 207
$ws.Range("E169").Value = 54
$ws.Range("H169").Value = 0

$ws.Range("A170").Value = "Lesoto"
$ws.Range("B170").Value = 256
$ws.Range("D170").Value = 48
$ws.Range("E170").Value = 205
$ws.Range("H170").Value = 3

$ws.Range("A171").Value = "Martinica"
$ws.Range("B171").Value = 255
$ws.Range("D171").Value = 98
$ws.Range("E171").Value = 142
$ws.Range("H171").Value = 15
